$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text fixups on existing rows -----------------------------------------

# RULE-100 details (E2): "Token" -> "Required token"
$ws.Range("E2").Value = "• Validation failures:`n• Required token 'apiId' not found in file: Properties\OCP\ITE.properties (case-sensitive: true)"

# RULE-101 details (E3): drop "test" from the expected values list
$ws.Range("E3").Value = "• Validation failures:`n• Property 'LogJsonFormat' found but value does not match expected values [true, false] in file: Properties\OCP\ITE.properties`n• Property 'anotherpropertycheck' found but value does not match expected values [somevalue] in file: Properties\OCP\ITE.properties"

# RULE-102 details (E4): "Token" -> "Required token" (x4)
$ws.Range("E4").Value = "• Validation failures:`n• Required token 'http.protocols=HTTPS' not found in file: Policies\TDV.policy (case-sensitive: true)`n• Required token 'http.private.port=8081' not found in file: Policies\TDV.policy (case-sensitive: true)`n• Required token 'http.protocols=HTTPS' not found in file: Policies\TDV1.policy (case-sensitive: true)`n• Required token 'http.private.port=8081' not found in file: Policies\TDV1.policy (case-sensitive: true)"

# --- Column B width: 42.18359375 -> 52.0 -----------------------------------
$ws.Columns.Item(2).ColumnWidth = 51.166666666666664

# --- New rows 9-12 -----------------------------------------------------------

$newRows = @(
    @{ Row = 9;  A = "RULE-107"; B = "Forbidden substring check for .properties files"; C = "HIGH"; D = "FAIL"; E = "• Validation failures:`n• Forbidden token 'fixmelater' found in file: Properties\OCP\ITE.properties (case-sensitive: true)" },
    @{ Row = 10; A = "RULE-108"; B = "Forbidden substring check for .policy files"; C = "HIGH"; D = "FAIL"; E = "• Validation failures:`n• Forbidden token 'deprecated.policy' found in file: Policies\TDV.policy (case-sensitive: true)" },
    @{ Row = 11; A = "RULE-109"; B = "Forbidden regex pattern (ip addresses) check in .properties files"; C = "HIGH"; D = "FAIL"; E = "• Forbidden token '^(?![\s]*[#!]).*\b(?:[0-9]{1,3}\.){3}[0-9]{1,3}\b' found in file: Properties\OCP\ITE.properties" },
    @{ Row = 12; A = "RULE-110"; B = "Forbidden regex pattern (ip addresses) check in .policy files"; C = "HIGH"; D = "FAIL"; E = "• Forbidden token '^(?![\s]*[#!]).*\b(?:[0-9]{1,3}\.){3}[0-9]{1,3}\b' found in file: Policies\TDV.policy" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range(("A{0}:E{0}" -f $row)).Interior.ColorIndex = 22
}
